$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Merge-Runs($shape) {
    $range = $shape.TextFrame.TextRange
    $text = $range.Text
    $len = $text.Length

    # Leading line/paragraph breaks show up as vertical-tab (chr 11)
    # characters inside TextRange.Text; keep them untouched and only
    # collapse the runs that make up the rest of the text.
    $skip = 0
    while ($skip -lt $len -and [int][char]$text[$skip] -eq 11) {
        $skip = $skip + 1
    }

    if ($skip -lt $len) {
        $rest = $range.Characters($skip + 1, $len - $skip)
        $rest.Text = $rest.Text
    }
}

# Title placeholder ("Testing" / " " / "custom" / " " / "properties")
# -> merge the split runs into a single run with the same text.
Merge-Runs($s.Shapes.Item(1))

# Subtitle placeholder ("A." / " " / "M." after two line breaks)
# -> merge the split runs into a single run with the same text,
#    keeping the leading line breaks untouched.
Merge-Runs($s.Shapes.Item(2))
